$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-13 22:52:21"
$wsZh.Range("H4").Value = "2016-03-13 22:52:53"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-13 22:52:24"
$wsDe.Range("H4").Value = "2016-03-13 22:53:00"
